# Fruta / hortaliza, semanal
# Insert a new weekly record at row 135 (pushing existing rows 135-147 down to
# 136-148) on the "Perejil - Feria Lagunitas de Puerto Montt" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 135..147 down by one to make room for the new record.
$ws.Rows(135).Insert()

# Populate the newly inserted row 135 with the new weekly entry.
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44449
$ws.Range("D135").NumberFormat = $ws.Range("D136").NumberFormat
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112044
$ws.Range("G135").Value = "Perejil"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 150
$ws.Range("K135").Value = 5000
$ws.Range("L135").Value = 5000
$ws.Range("M135").Value = 5000
$ws.Range("N135").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O135").Value = "Región Metropolitana"
$ws.Range("P135").Value = 1667
$ws.Range("Q135").Value = 3
$ws.Range("R135").Value = "Hortaliza"
